$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2 - convert comma-decimal strings into real numeric values
$ws.Range("C2").Value = 935.9400000000001
$ws.Range("E2").Value = 91.38
$ws.Range("G2").Value = 25
$ws.Range("H2").Value = 160.34
$ws.Range("I2").Value = 211.33
$ws.Range("J2").Value = 243.24
$ws.Range("K2").Value = 139.95
$ws.Range("L2").Value = 785.37
$ws.Range("M2").Value = 45912.82603366111

# Row 3 - convert comma-decimal strings into real numeric values
$ws.Range("C3").Value = 426.97
$ws.Range("E3").Value = 33.64
$ws.Range("G3").Value = 25
$ws.Range("H3").Value = 54.14
$ws.Range("I3").Value = 192.16

# J3 / L3 become blank text cells (content removed, cell itself retained)
$ws.Range("J3").Value = "'"
$ws.Range("J3").Style = "Normal"
$ws.Range("K3").Value = 284.81
$ws.Range("L3").Value = "'"
$ws.Range("L3").Style = "Normal"

$ws.Range("M3").Value = 45912.82603366111
